#
# Populate the LoginPage sheet with the new locator rows, switch the
# "active" sheet/selection from AddEmployeePage back to LoginPage, and
# refresh a couple of leftover selections, matching the author's
# "Add Employee page also in progress" commit.
#

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LoginPage")
$ws2 = $wb.Worksheets.Item("AddEmployeePage")

# ---------------------------------------------------------------------
# LoginPage (sheet1): rewrite column A (locator names) and fill in the
# new id/cssselector/xpath/classname/tagname columns for every locator.
# ---------------------------------------------------------------------

$ws1.Range("A2").Value = "txt_username"
$ws1.Range("C2").Value = "username"
$ws1.Range("D2").Value = "input[placeholder='username']"
$ws1.Range("E2").Value = "//input[@placeholder='username']"
$ws1.Range("F2").Value = "oxd-input oxd-input--active"
$ws1.Range("G2").Value = "input"

$ws1.Range("A3").Value = "txt_password"
$ws1.Range("C3").Value = "password"
$ws1.Range("D3").Value = "input[placeholder='password']"
$ws1.Range("E3").Value = "//input[@placeholder='password']"
$ws1.Range("F3").Value = "oxd-input oxd-input--active"
$ws1.Range("G3").Value = "input"

$ws1.Range("A4").Value = "login_btn"
$ws1.Range("D4").Value = "button[type='submit']"
$ws1.Range("E4").Value = "//button[@type='submit']"
$ws1.Range("F4").Value = "oxd-button oxd-button--medium oxd-button--main orangehrm-login-button"
$ws1.Range("G4").Value = "button"

$ws1.Range("A5").Value = "page_heading"
$ws1.Range("D5").Value = ".oxd-text.oxd-text--h5.orangehrm-login-title"
$ws1.Range("E5").Value = "//h5[text()='Login']"
$ws1.Range("F5").Value = "oxd-text oxd-text--h5 orangehrm-login-title"
$ws1.Range("G5").Value = "h5"

$ws1.Range("A6").Value = "invalid_txt"
$ws1.Range("D6").Value = ".oxd-text.oxd-text--p.oxd-alert-content-text"
$ws1.Range("E6").Value = "//*[text()='Invalid credentials']"
$ws1.Range("F6").Value = "oxd-text oxd-text--p oxd-alert-content-text"
$ws1.Range("G6").Value = "p"

$ws1.Range("A7").Value = "required_msg"
$ws1.Range("D7").Value = ".oxd-text.oxd-text--span.oxd-input-field-error-message.oxd-input-group__message"
$ws1.Range("E7").Value = "//*[text()='Required']"
$ws1.Range("F7").Value = "oxd-text oxd-text--span oxd-input-field-error-message oxd-input-group__message"
$ws1.Range("G7").Value = "span"

# Wrap the long selector/xpath/classname text, same as the header and
# the AddEmployeePage sheet already do.
$ws1.Range("D1:F7").WrapText = $true

# Row heights so the wrapped text is fully visible (matches the sizes
# already used on AddEmployeePage for equivalent content).
$ws1.Rows.Item(2).RowHeight = 28.8
$ws1.Rows.Item(3).RowHeight = 28.8
$ws1.Rows.Item(4).RowHeight = 43.2
$ws1.Rows.Item(5).RowHeight = 28.8
$ws1.Rows.Item(6).RowHeight = 28.8
$ws1.Rows.Item(7).RowHeight = 57.6

# Widen the cssselector/xpath/classname columns to fit the new content.
$ws1.Columns.Item(4).ColumnWidth = 25.66406
$ws1.Columns.Item(5).ColumnWidth = 29.21875
$ws1.Columns.Item(6).ColumnWidth = 23.21875

# ---------------------------------------------------------------------
# Fix up sheet selections / active sheet: LoginPage becomes the active
# tab again (with D7 selected), AddEmployeePage keeps its own last
# selection (E6) but is no longer the active tab.
# ---------------------------------------------------------------------

$ws2.Activate() | Out-Null
$ws2.Range("E6").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("D7").Select() | Out-Null

Write-Output "LoginPage locators populated; AddEmployeePage selection updated."
